$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename sheet to reflect new "through" date
$ws.Name = "Through 2022-04-14"

# Update the header label for the running-total column (shared string used in I1)
$ws.Range("I1").Value = "2022 (through 04-14)"

# Update April's running total (row 5) and the overall Total row (row 14)
$ws.Range("I5").Value = 59
$ws.Range("I14").Value = 493
